# GPLIM-2588 Fix spreadsheet headers.
#
# Rename the two header cells that drive manifest-upload column matching:
#   A1: "Sample ID" -> "Specimen_Number"
#   F1: "T/N"       -> "SAMPLE_TYPE"
#
# We stage each new header text in a scratch cell, Copy it, and
# PasteSpecial(Values) it into the target header cell. This swaps the
# shared-string reference for the cell while leaving the destination
# cell's existing formatting/style untouched (a plain `.Value =`
# assignment would re-derive the cell's style from the new text and
# split off a new style record, which is not what happened here).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buick Example")

$scratch = $ws.Range("H1")

$scratch.Value = "Specimen_Number"
$scratch.Copy()
$ws.Range("A1").PasteSpecial(-4163)

$scratch.Value = "SAMPLE_TYPE"
$scratch.Copy()
$ws.Range("F1").PasteSpecial(-4163)

$scratch.ClearContents()

# Matches the author's final on-screen selection after the edit.
$ws.Range("F2").Select()
